# Auto-generated edit script updating FFXIV Excalibur profit calculations
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4195.4736
$ws.Range("I138").Value = 2524.2222
$ws.Range("J138").Value = 5699.6
$ws.Range("K138").Value = 7572.6666
$ws.Range("L138").Value = 17098.8
$ws.Range("M138").Value = -2432.6666
$ws.Range("N138").Value = -27378.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 232.27272
$ws.Range("I4").Value = 266.44446
$ws.Range("J4").Value = 78.5
$ws.Range("K4").Value = 266.44446
$ws.Range("L4").Value = 78.5
$ws.Range("M4").Value = -150.44446
$ws.Range("N4").Value = -310.5
$ws.Range("H36").Value = 7000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 7000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 7000
$ws.Range("M36").Value = $null
$ws.Range("N36").Value = -7692
$ws.Range("H45").Value = 3635.3635
$ws.Range("I45").Value = 2462.8333
$ws.Range("J45").Value = 5042.4
$ws.Range("K45").Value = 2462.8333
$ws.Range("L45").Value = 5042.4
$ws.Range("M45").Value = -2085.8333
$ws.Range("N45").Value = -5796.4
$ws.Range("H61").Value = 11116776
$ws.Range("I61").Value = 16670164
$ws.Range("J61").Value = 9999
$ws.Range("K61").Value = 16670164
$ws.Range("L61").Value = 9999
$ws.Range("M61").Value = -16669952
$ws.Range("N61").Value = -10423
$ws.Range("H63").Value = 5682.1113
$ws.Range("I63").Value = 2942.3333
$ws.Range("J63").Value = 8421.888999999999
$ws.Range("K63").Value = 2942.3333
$ws.Range("L63").Value = 8421.888999999999
$ws.Range("M63").Value = -2256.3333
$ws.Range("N63").Value = -9793.888999999999
$ws.Range("H66").Value = 5682.1113
$ws.Range("I66").Value = 2942.3333
$ws.Range("J66").Value = 8421.888999999999
$ws.Range("K66").Value = 14711.6665
$ws.Range("L66").Value = 42109.44499999999
$ws.Range("M66").Value = -11279.6665
$ws.Range("N66").Value = -48973.44499999999
$ws.Range("H74").Value = 2576.925
$ws.Range("I74").Value = 1299.3846
$ws.Range("J74").Value = 4949.5
$ws.Range("K74").Value = 1299.3846
$ws.Range("L74").Value = 4949.5
$ws.Range("M74").Value = -425.3846000000001
$ws.Range("N74").Value = -6697.5
$ws.Range("H77").Value = 2576.925
$ws.Range("I77").Value = 1299.3846
$ws.Range("J77").Value = 4949.5
$ws.Range("K77").Value = 6496.923000000001
$ws.Range("L77").Value = 24747.5
$ws.Range("M77").Value = -2128.923000000001
$ws.Range("N77").Value = -33483.5
$ws.Range("H97").Value = 1303
$ws.Range("I97").Value = 1303
$ws.Range("K97").Value = 1303
$ws.Range("M97").Value = -807
$ws.Range("H122").Value = 3295.68
$ws.Range("I122").Value = 2523.7058
$ws.Range("J122").Value = 4936.125
$ws.Range("K122").Value = 7571.117400000001
$ws.Range("L122").Value = 14808.375
$ws.Range("M122").Value = -5121.117400000001
$ws.Range("N122").Value = -19708.375
$ws.Range("H132").Value = 691376.5
$ws.Range("I132").Value = 715925.7
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 2147777.1
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -2145247.1
$ws.Range("N132").Value = -17060
$ws.Range("H136").Value = 11116776
$ws.Range("I136").Value = 16670164
$ws.Range("J136").Value = 9999
$ws.Range("K136").Value = 50010492
$ws.Range("L136").Value = 29997
$ws.Range("M136").Value = -50007942
$ws.Range("N136").Value = -35097

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2083.5
$ws.Range("I94").Value = 2055.5
$ws.Range("J94").Value = 2097.5
$ws.Range("K94").Value = 2055.5
$ws.Range("L94").Value = 2097.5
$ws.Range("M94").Value = -1604.5
$ws.Range("N94").Value = -2999.5
$ws.Range("H107").Value = 2734.2144
$ws.Range("I107").Value = 2759.923
$ws.Range("J107").Value = 2400
$ws.Range("K107").Value = 2759.923
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = -839.9229999999998
$ws.Range("N107").Value = -6240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1201.5625
$ws.Range("J16").Value = 1708
$ws.Range("L16").Value = 1708
$ws.Range("N16").Value = -2282
$ws.Range("H31").Value = 9059.884
$ws.Range("I31").Value = 3726.9736
$ws.Range("J31").Value = 18271.273
$ws.Range("K31").Value = 3726.9736
$ws.Range("L31").Value = 18271.273
$ws.Range("M31").Value = -3431.9736
$ws.Range("N31").Value = -18861.273
$ws.Range("H34").Value = 9059.884
$ws.Range("I34").Value = 3726.9736
$ws.Range("J34").Value = 18271.273
$ws.Range("K34").Value = 3726.9736
$ws.Range("L34").Value = 18271.273
$ws.Range("M34").Value = -3524.9736
$ws.Range("N34").Value = -18675.273
$ws.Range("H58").Value = 443117.78
$ws.Range("I58").Value = 727751.7
$ws.Range("K58").Value = 727751.7
$ws.Range("M58").Value = -727548.7
$ws.Range("H107").Value = 700.5185
$ws.Range("I107").Value = 575.4706
$ws.Range("J107").Value = 913.1
$ws.Range("K107").Value = 575.4706
$ws.Range("L107").Value = 913.1
$ws.Range("M107").Value = 1344.5294
$ws.Range("N107").Value = -4753.1
$ws.Range("H113").Value = 1201.5625
$ws.Range("J113").Value = 1708
$ws.Range("L113").Value = 1708
$ws.Range("N113").Value = -6048
$ws.Range("H122").Value = 2457.88
$ws.Range("I122").Value = 1503.5625
$ws.Range("J122").Value = 4154.4443
$ws.Range("K122").Value = 4510.6875
$ws.Range("L122").Value = 12463.3329
$ws.Range("M122").Value = -2060.6875
$ws.Range("N122").Value = -17363.3329
$ws.Range("H132").Value = 19125.715
$ws.Range("I132").Value = 1751.4
$ws.Range("J132").Value = 62561.5
$ws.Range("K132").Value = 5254.200000000001
$ws.Range("L132").Value = 187684.5
$ws.Range("M132").Value = -2724.200000000001
$ws.Range("N132").Value = -192744.5
$ws.Range("H134").Value = 7138.727
$ws.Range("I134").Value = 7138.727
$ws.Range("K134").Value = 21416.181
$ws.Range("M134").Value = -18881.181
$ws.Range("H136").Value = 443117.78
$ws.Range("I136").Value = 727751.7
$ws.Range("K136").Value = 2183255.1
$ws.Range("M136").Value = -2180705.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 6388.763
$ws.Range("I68").Value = 2090.6667
$ws.Range("K68").Value = 6272.000100000001
$ws.Range("M68").Value = -5461.000100000001
$ws.Range("H71").Value = 6388.763
$ws.Range("I71").Value = 2090.6667
$ws.Range("K71").Value = 18816.0003
$ws.Range("M71").Value = -14760.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 37230.2
$ws.Range("J45").Value = 37230.2
$ws.Range("L45").Value = 37230.2
$ws.Range("N45").Value = -38348.2
$ws.Range("H107").Value = 22147.375
$ws.Range("I107").Value = 29407.084
$ws.Range("J107").Value = 368.25
$ws.Range("K107").Value = 29407.084
$ws.Range("L107").Value = 368.25
$ws.Range("M107").Value = -27487.084
$ws.Range("N107").Value = -4208.25
$ws.Range("H132").Value = 432663.2
$ws.Range("I132").Value = 448613.16
$ws.Range("J132").Value = 2014
$ws.Range("K132").Value = 1345839.48
$ws.Range("L132").Value = 6042
$ws.Range("M132").Value = -1343309.48
$ws.Range("N132").Value = -11102

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1058116.1
$ws.Range("I132").Value = 1452243.2
$ws.Range("K132").Value = 4356729.6
$ws.Range("M132").Value = -4354199.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1713.4445
$ws.Range("I107").Value = 923.4666999999999
$ws.Range("J107").Value = 5663.3335
$ws.Range("K107").Value = 2770.4001
$ws.Range("L107").Value = 16990.0005
$ws.Range("M107").Value = -850.4000999999998
$ws.Range("N107").Value = -20830.0005
$ws.Range("H132").Value = 4195309.5
$ws.Range("I132").Value = 4576088
$ws.Range("J132").Value = 6747.25
$ws.Range("K132").Value = 13728264
$ws.Range("L132").Value = 20241.75
$ws.Range("M132").Value = -13725734
$ws.Range("N132").Value = -25301.75

